$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first worker's (KARINA MARTINEZ AGAMEZ / 33104063) 13 rows of
# period data entirely - the new workbook only carries ARIEL ENRIQUE CASTRO
# VEGA's account statement. Deleting rows 16:28 shifts the remaining
# ARIEL data (previously rows 29:41) up into rows 16:28, and also shifts
# the two footer rows (previously 46:47) up into rows 33:34 automatically,
# along with the sheet dimension and merged cell ranges.
$ws.Rows("16:28").Delete()

# Updated summary figures for the single remaining worker.
$ws.Range("C13").Value = 1
$ws.Range("E11").Value = 939018

# Rewrite the 13 period rows (now 16:28) in chronological order with the
# refreshed Valor Mora / Salario Basico figures.
$periods = @("1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806")
$valorMora = @(80000,80000,80000,80000,80000,80000,80000,80000,80000,80000,80000,29509,29509)
$salarioBasico = 877803

for ($i = 0; $i -lt 13; $i++) {
    $r = 16 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "12639272"
    $ws.Range("D$r").Value = "ARIEL ENRIQUE CASTRO VEGA"
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("F$r").Value = $valorMora[$i]
    $ws.Range("G$r").Value = $salarioBasico
}
